$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$c = $t.Cell(1,1)
$c.Range.Find.Execute("4+44=", $true, $false, $false, $false, $false, $true, 0, $false, "88+8=", 1) | Out-Null
$c = $t.Cell(1,2)
$c.Range.Find.Execute("68-2=", $true, $false, $false, $false, $false, $true, 0, $false, "52+13=", 1) | Out-Null
$c = $t.Cell(1,3)
$c.Range.Find.Execute("21+55=", $true, $false, $false, $false, $false, $true, 0, $false, "17+58=", 1) | Out-Null
$c = $t.Cell(1,4)
$c.Range.Find.Execute("17+29=", $true, $false, $false, $false, $false, $true, 0, $false, "42-14=", 1) | Out-Null
$c = $t.Cell(1,5)
$c.Range.Find.Execute("82-0=", $true, $false, $false, $false, $false, $true, 0, $false, "41+8=", 1) | Out-Null
$c = $t.Cell(2,1)
$c.Range.Find.Execute("29+30=", $true, $false, $false, $false, $false, $true, 0, $false, "51-12=", 1) | Out-Null
$c = $t.Cell(2,2)
$c.Range.Find.Execute("19+36=", $true, $false, $false, $false, $false, $true, 0, $false, "12+11=", 1) | Out-Null
$c = $t.Cell(2,3)
$c.Range.Find.Execute("36-8=", $true, $false, $false, $false, $false, $true, 0, $false, "97-38=", 1) | Out-Null
$c = $t.Cell(2,4)
$c.Range.Find.Execute("54-24=", $true, $false, $false, $false, $false, $true, 0, $false, "65+34=", 1) | Out-Null
$c = $t.Cell(2,5)
$c.Range.Find.Execute("42+17=", $true, $false, $false, $false, $false, $true, 0, $false, "69-28=", 1) | Out-Null
$c = $t.Cell(3,1)
$c.Range.Find.Execute("42+1=", $true, $false, $false, $false, $false, $true, 0, $false, "77+4=", 1) | Out-Null
$c = $t.Cell(3,2)
$c.Range.Find.Execute("28+52=", $true, $false, $false, $false, $false, $true, 0, $false, "35-25=", 1) | Out-Null
$c = $t.Cell(3,3)
$c.Range.Find.Execute("13+22=", $true, $false, $false, $false, $false, $true, 0, $false, "40-17=", 1) | Out-Null
$c = $t.Cell(3,4)
$c.Range.Find.Execute("41+3=", $true, $false, $false, $false, $false, $true, 0, $false, "14+46=", 1) | Out-Null
$c = $t.Cell(3,5)
$c.Range.Find.Execute("11+7=", $true, $false, $false, $false, $false, $true, 0, $false, "66+18=", 1) | Out-Null
$c = $t.Cell(4,1)
$c.Range.Find.Execute("91-11=", $true, $false, $false, $false, $false, $true, 0, $false, "43+32=", 1) | Out-Null
$c = $t.Cell(4,2)
$c.Range.Find.Execute("67-62=", $true, $false, $false, $false, $false, $true, 0, $false, "83-31=", 1) | Out-Null
$c = $t.Cell(4,3)
$c.Range.Find.Execute("39-30=", $true, $false, $false, $false, $false, $true, 0, $false, "71-0=", 1) | Out-Null
$c = $t.Cell(4,4)
$c.Range.Find.Execute("85-59=", $true, $false, $false, $false, $false, $true, 0, $false, "65-54=", 1) | Out-Null
$c = $t.Cell(4,5)
$c.Range.Find.Execute("26+29=", $true, $false, $false, $false, $false, $true, 0, $false, "23+1=", 1) | Out-Null
$c = $t.Cell(5,1)
$c.Range.Find.Execute("54-27=", $true, $false, $false, $false, $false, $true, 0, $false, "90-90=", 1) | Out-Null
$c = $t.Cell(5,2)
$c.Range.Find.Execute("47+18=", $true, $false, $false, $false, $false, $true, 0, $false, "65+27=", 1) | Out-Null
$c = $t.Cell(5,3)
$c.Range.Find.Execute("24+28=", $true, $false, $false, $false, $false, $true, 0, $false, "77+2=", 1) | Out-Null
$c = $t.Cell(5,4)
$c.Range.Find.Execute("28+51=", $true, $false, $false, $false, $false, $true, 0, $false, "5+41=", 1) | Out-Null
$c = $t.Cell(5,5)
$c.Range.Find.Execute("81-56=", $true, $false, $false, $false, $false, $true, 0, $false, "52+40=", 1) | Out-Null
$c = $t.Cell(6,1)
$c.Range.Find.Execute("72-36=", $true, $false, $false, $false, $false, $true, 0, $false, "58-38=", 1) | Out-Null
$c = $t.Cell(6,2)
$c.Range.Find.Execute("84-8=", $true, $false, $false, $false, $false, $true, 0, $false, "53-8=", 1) | Out-Null
$c = $t.Cell(6,3)
$c.Range.Find.Execute("17+28=", $true, $false, $false, $false, $false, $true, 0, $false, "83-77=", 1) | Out-Null
$c = $t.Cell(6,4)
$c.Range.Find.Execute("16+56=", $true, $false, $false, $false, $false, $true, 0, $false, "97-49=", 1) | Out-Null
$c = $t.Cell(6,5)
$c.Range.Find.Execute("67-0=", $true, $false, $false, $false, $false, $true, 0, $false, "95-92=", 1) | Out-Null
$c = $t.Cell(7,1)
$c.Range.Find.Execute("41+2=", $true, $false, $false, $false, $false, $true, 0, $false, "5+76=", 1) | Out-Null
$c = $t.Cell(7,2)
$c.Range.Find.Execute("83-43=", $true, $false, $false, $false, $false, $true, 0, $false, "83-36=", 1) | Out-Null
$c = $t.Cell(7,3)
$c.Range.Find.Execute("96-57=", $true, $false, $false, $false, $false, $true, 0, $false, "30+43=", 1) | Out-Null
$c = $t.Cell(7,4)
$c.Range.Find.Execute("12-3=", $true, $false, $false, $false, $false, $true, 0, $false, "72-43=", 1) | Out-Null
$c = $t.Cell(7,5)
$c.Range.Find.Execute("97-50=", $true, $false, $false, $false, $false, $true, 0, $false, "6+52=", 1) | Out-Null
$c = $t.Cell(8,1)
$c.Range.Find.Execute("1+46=", $true, $false, $false, $false, $false, $true, 0, $false, "38+54=", 1) | Out-Null
$c = $t.Cell(8,2)
$c.Range.Find.Execute("73-39=", $true, $false, $false, $false, $false, $true, 0, $false, "39+4=", 1) | Out-Null
$c = $t.Cell(8,3)
$c.Range.Find.Execute("41+39=", $true, $false, $false, $false, $false, $true, 0, $false, "9+44=", 1) | Out-Null
$c = $t.Cell(8,4)
$c.Range.Find.Execute("97-47=", $true, $false, $false, $false, $false, $true, 0, $false, "95-52=", 1) | Out-Null
$c = $t.Cell(8,5)
$c.Range.Find.Execute("92-69=", $true, $false, $false, $false, $false, $true, 0, $false, "15+25=", 1) | Out-Null
$c = $t.Cell(9,1)
$c.Range.Find.Execute("2+9=", $true, $false, $false, $false, $false, $true, 0, $false, "72-0=", 1) | Out-Null
$c = $t.Cell(9,2)
$c.Range.Find.Execute("7+19=", $true, $false, $false, $false, $false, $true, 0, $false, "55-35=", 1) | Out-Null
$c = $t.Cell(9,3)
$c.Range.Find.Execute("64-1=", $true, $false, $false, $false, $false, $true, 0, $false, "63-51=", 1) | Out-Null
$c = $t.Cell(9,4)
$c.Range.Find.Execute("9+47=", $true, $false, $false, $false, $false, $true, 0, $false, "29+8=", 1) | Out-Null
$c = $t.Cell(9,5)
$c.Range.Find.Execute("1+72=", $true, $false, $false, $false, $false, $true, 0, $false, "96-35=", 1) | Out-Null
$c = $t.Cell(10,1)
$c.Range.Find.Execute("8+20=", $true, $false, $false, $false, $false, $true, 0, $false, "96-67=", 1) | Out-Null
$c = $t.Cell(10,2)
$c.Range.Find.Execute("73-66=", $true, $false, $false, $false, $false, $true, 0, $false, "36-30=", 1) | Out-Null
$c = $t.Cell(10,3)
$c.Range.Find.Execute("79-74=", $true, $false, $false, $false, $false, $true, 0, $false, "14+28=", 1) | Out-Null
$c = $t.Cell(10,4)
$c.Range.Find.Execute("17+7=", $true, $false, $false, $false, $false, $true, 0, $false, "90-24=", 1) | Out-Null
$c = $t.Cell(10,5)
$c.Range.Find.Execute("4+68=", $true, $false, $false, $false, $false, $true, 0, $false, "2+21=", 1) | Out-Null
$c = $t.Cell(11,1)
$c.Range.Find.Execute("67-67=", $true, $false, $false, $false, $false, $true, 0, $false, "76+23=", 1) | Out-Null
$c = $t.Cell(11,2)
$c.Range.Find.Execute("4+57=", $true, $false, $false, $false, $false, $true, 0, $false, "31+47=", 1) | Out-Null
$c = $t.Cell(11,3)
$c.Range.Find.Execute("79-9=", $true, $false, $false, $false, $false, $true, 0, $false, "92-67=", 1) | Out-Null
$c = $t.Cell(11,4)
$c.Range.Find.Execute("64+2=", $true, $false, $false, $false, $false, $true, 0, $false, "10+64=", 1) | Out-Null
$c = $t.Cell(11,5)
$c.Range.Find.Execute("14+83=", $true, $false, $false, $false, $false, $true, 0, $false, "6+84=", 1) | Out-Null
$c = $t.Cell(12,1)
$c.Range.Find.Execute("62-50=", $true, $false, $false, $false, $false, $true, 0, $false, "46+2=", 1) | Out-Null
$c = $t.Cell(12,2)
$c.Range.Find.Execute("82-27=", $true, $false, $false, $false, $false, $true, 0, $false, "63+20=", 1) | Out-Null
$c = $t.Cell(12,3)
$c.Range.Find.Execute("52+46=", $true, $false, $false, $false, $false, $true, 0, $false, "47+35=", 1) | Out-Null
$c = $t.Cell(12,4)
$c.Range.Find.Execute("38-37=", $true, $false, $false, $false, $false, $true, 0, $false, "2+88=", 1) | Out-Null
$c = $t.Cell(12,5)
$c.Range.Find.Execute("88-52=", $true, $false, $false, $false, $false, $true, 0, $false, "30+29=", 1) | Out-Null
$c = $t.Cell(13,1)
$c.Range.Find.Execute("31-5=", $true, $false, $false, $false, $false, $true, 0, $false, "73-12=", 1) | Out-Null
$c = $t.Cell(13,2)
$c.Range.Find.Execute("67+31=", $true, $false, $false, $false, $false, $true, 0, $false, "46+34=", 1) | Out-Null
$c = $t.Cell(13,3)
$c.Range.Find.Execute("67-49=", $true, $false, $false, $false, $false, $true, 0, $false, "63-22=", 1) | Out-Null
$c = $t.Cell(13,4)
$c.Range.Find.Execute("14+2=", $true, $false, $false, $false, $false, $true, 0, $false, "47-23=", 1) | Out-Null
$c = $t.Cell(13,5)
$c.Range.Find.Execute("75-64=", $true, $false, $false, $false, $false, $true, 0, $false, "50-24=", 1) | Out-Null
$c = $t.Cell(14,1)
$c.Range.Find.Execute("64+13=", $true, $false, $false, $false, $false, $true, 0, $false, "35-22=", 1) | Out-Null
$c = $t.Cell(14,2)
$c.Range.Find.Execute("22+46=", $true, $false, $false, $false, $false, $true, 0, $false, "85-20=", 1) | Out-Null
$c = $t.Cell(14,3)
$c.Range.Find.Execute("49+1=", $true, $false, $false, $false, $false, $true, 0, $false, "34+42=", 1) | Out-Null
$c = $t.Cell(14,4)
$c.Range.Find.Execute("97-46=", $true, $false, $false, $false, $false, $true, 0, $false, "95-49=", 1) | Out-Null
$c = $t.Cell(14,5)
$c.Range.Find.Execute("27-8=", $true, $false, $false, $false, $false, $true, 0, $false, "75-23=", 1) | Out-Null
$c = $t.Cell(15,1)
$c.Range.Find.Execute("27+40=", $true, $false, $false, $false, $false, $true, 0, $false, "65-14=", 1) | Out-Null
$c = $t.Cell(15,2)
$c.Range.Find.Execute("21+21=", $true, $false, $false, $false, $false, $true, 0, $false, "90-20=", 1) | Out-Null
$c = $t.Cell(15,3)
$c.Range.Find.Execute("61-7=", $true, $false, $false, $false, $false, $true, 0, $false, "95-90=", 1) | Out-Null
$c = $t.Cell(15,4)
$c.Range.Find.Execute("86-75=", $true, $false, $false, $false, $false, $true, 0, $false, "15+44=", 1) | Out-Null
$c = $t.Cell(15,5)
$c.Range.Find.Execute("79-35=", $true, $false, $false, $false, $false, $true, 0, $false, "40+54=", 1) | Out-Null
$c = $t.Cell(16,1)
$c.Range.Find.Execute("16+40=", $true, $false, $false, $false, $false, $true, 0, $false, "59-5=", 1) | Out-Null
$c = $t.Cell(16,2)
$c.Range.Find.Execute("96-69=", $true, $false, $false, $false, $false, $true, 0, $false, "36-6=", 1) | Out-Null
$c = $t.Cell(16,3)
$c.Range.Find.Execute("55-25=", $true, $false, $false, $false, $false, $true, 0, $false, "16+70=", 1) | Out-Null
$c = $t.Cell(16,4)
$c.Range.Find.Execute("82-64=", $true, $false, $false, $false, $false, $true, 0, $false, "73-58=", 1) | Out-Null
$c = $t.Cell(16,5)
$c.Range.Find.Execute("77-12=", $true, $false, $false, $false, $false, $true, 0, $false, "19+41=", 1) | Out-Null
$c = $t.Cell(17,1)
$c.Range.Find.Execute("89-58=", $true, $false, $false, $false, $false, $true, 0, $false, "87-14=", 1) | Out-Null
$c = $t.Cell(17,2)
$c.Range.Find.Execute("20+25=", $true, $false, $false, $false, $false, $true, 0, $false, "77-62=", 1) | Out-Null
$c = $t.Cell(17,3)
$c.Range.Find.Execute("85-85=", $true, $false, $false, $false, $false, $true, 0, $false, "1+60=", 1) | Out-Null
$c = $t.Cell(17,4)
$c.Range.Find.Execute("52-1=", $true, $false, $false, $false, $false, $true, 0, $false, "96-10=", 1) | Out-Null
$c = $t.Cell(17,5)
$c.Range.Find.Execute("47-1=", $true, $false, $false, $false, $false, $true, 0, $false, "55-35=", 1) | Out-Null
$c = $t.Cell(18,1)
$c.Range.Find.Execute("69-13=", $true, $false, $false, $false, $false, $true, 0, $false, "3+74=", 1) | Out-Null
$c = $t.Cell(18,2)
$c.Range.Find.Execute("9+19=", $true, $false, $false, $false, $false, $true, 0, $false, "92-44=", 1) | Out-Null
$c = $t.Cell(18,3)
$c.Range.Find.Execute("47-37=", $true, $false, $false, $false, $false, $true, 0, $false, "51+19=", 1) | Out-Null
$c = $t.Cell(18,4)
$c.Range.Find.Execute("94-4=", $true, $false, $false, $false, $false, $true, 0, $false, "77-27=", 1) | Out-Null
$c = $t.Cell(18,5)
$c.Range.Find.Execute("56-50=", $true, $false, $false, $false, $false, $true, 0, $false, "84+12=", 1) | Out-Null
$c = $t.Cell(19,1)
$c.Range.Find.Execute("74+12=", $true, $false, $false, $false, $false, $true, 0, $false, "18-8=", 1) | Out-Null
$c = $t.Cell(19,2)
$c.Range.Find.Execute("17+19=", $true, $false, $false, $false, $false, $true, 0, $false, "57+2=", 1) | Out-Null
$c = $t.Cell(19,3)
$c.Range.Find.Execute("66+3=", $true, $false, $false, $false, $false, $true, 0, $false, "13+77=", 1) | Out-Null
$c = $t.Cell(19,4)
$c.Range.Find.Execute("29+4=", $true, $false, $false, $false, $false, $true, 0, $false, "96-54=", 1) | Out-Null
$c = $t.Cell(19,5)
$c.Range.Find.Execute("89-28=", $true, $false, $false, $false, $false, $true, 0, $false, "95-28=", 1) | Out-Null
$c = $t.Cell(20,1)
$c.Range.Find.Execute("97-33=", $true, $false, $false, $false, $false, $true, 0, $false, "47+43=", 1) | Out-Null
$c = $t.Cell(20,2)
$c.Range.Find.Execute("61+19=", $true, $false, $false, $false, $false, $true, 0, $false, "10-4=", 1) | Out-Null
$c = $t.Cell(20,3)
$c.Range.Find.Execute("47+38=", $true, $false, $false, $false, $false, $true, 0, $false, "96-62=", 1) | Out-Null
$c = $t.Cell(20,4)
$c.Range.Find.Execute("99-1=", $true, $false, $false, $false, $false, $true, 0, $false, "17+15=", 1) | Out-Null
$c = $t.Cell(20,5)
$c.Range.Find.Execute("96-11=", $true, $false, $false, $false, $false, $true, 0, $false, "86-8=", 1) | Out-Null
